$d = $word.ActiveDocument

function Replace-Exact($findText, $newText) {
    $r = $d.Content
    $found = $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $r.Text = $newText
    }
}

Replace-Exact "joe@gmail.com" "bobrossfake@gmail.com"
Replace-Exact "Password: joepassword" "Password: bobbyword"
Replace-Exact "Firstname: joe" "Firstname: bob"
Replace-Exact "Surname: munkey" "Surname: ross"
Replace-Exact "Address: 57 road" "Address: 12 angel street"
Replace-Exact "Postcode: bt35 rf4" "Postcode: yr67 1er"
Replace-Exact "Age: 26" "Age: 33"
Replace-Exact "Group: 6" "Group: 7"
